# Update countries & provincias Spain
#
# 1) Swap the Suecia/Rusia rows (22 <-> 23) so "Rusia" now sorts before
#    "Suecia" in the shared-strings table / row order, and refresh Rusia's
#    daily figures to the newer snapshot.
# 2) Bump the "Datos actualizados" timestamp from 09:22 to 09:52.
# 3) Refresh a handful of per-country case counts that were updated in the
#    same data pull (Moldavia, Barein, Bosnia y Herzegovina, Letonia,
#    Sri Lanka).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row 22 becomes Rusia (with refreshed figures), row 23 becomes Suecia
#        (with the figures Rusia used to have before the refresh). ---
$ws.Range("A22").Value = "Rusia"
$ws.Range("B22").Value = 8672
$ws.Range("C22").Value = 1175
$ws.Range("D22").Value = 580
$ws.Range("E22").Value = 8029
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 63

$ws.Range("A23").Value = "Suecia"
$ws.Range("B23").Value = 7693
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 205
$ws.Range("E23").Value = 6897
$ws.Range("F23").Value = 640
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 591

# --- 2) Timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 09:52"

# --- 3) Other country figure refreshes ---

# Row 65: Moldavia
$ws.Range("E65").Value = 992
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 24

# Row 71: Barein
$ws.Range("D71").Value = 465
$ws.Range("E71").Value = 341
$ws.Range("F71").Value = 3

# Row 72: Bosnia y Herzegovina
$ws.Range("B72").Value = 777
$ws.Range("C72").Value = 13
$ws.Range("D72").Value = 77
$ws.Range("E72").Value = 667

# Row 82: Letonia
$ws.Range("B82").Value = 577
$ws.Range("C82").Value = 29
$ws.Range("E82").Value = 559

# Row 112: Sri Lanka
$ws.Range("B112").Value = 186
$ws.Range("C112").Value = 1
$ws.Range("E112").Value = 138
